$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.043.77'
$ws.Range("E2").Value = '  -0.67%  '

$ws.Range("D3").Value = '1.822.65'
$ws.Range("E3").Value = '  +2.61%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '337.27'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").Value = '0.4244'
$ws.Range("E7").Value = '  +10.90%  '

$ws.Range("D8").Value = '0.3529'
$ws.Range("E8").Value = '  +3.02%  '

$ws.Range("D9").Value = '45.62'
$ws.Range("E9").Value = '  -2.71%  '

$ws.Range("D10").Value = '1.158'
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("D11").Value = '0.07510'
$ws.Range("E11").Value = '  +1.57%  '

$ws.Range("D12").Value = '23.05'
$ws.Range("E12").Value = '  -2.38%  '

$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").Value = '6.308'
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").Value = '7.319'
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").Value = '1.819.16'
$ws.Range("E16").Value = '  +2.30%  '

$ws.Range("E17").Value = '  +1.48%  '

$ws.Range("D18").Value = '0.06712'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("D19").Value = '82.80'
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").Value = '17.47'
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").Value = '6.413'
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").Value = '28.109.52'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("D24").Value = '11.92'
$ws.Range("E24").Value = '  -1.51%  '

$ws.Range("D25").Value = '2.406'
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("D26").Value = '2.505'
$ws.Range("E26").Value = '  +3.99%  '

$ws.Range("D27").Value = '20.95'
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D28").Value = '156.34'
$ws.Range("E28").Value = '  +1.83%  '

$ws.Range("D29").Value = '2.027.89'
$ws.Range("E29").Value = '  +2.49%  '

$ws.Range("D30").Value = '1.316'
$ws.Range("E30").Value = '  -7.86%  '

$ws.Range("D31").Value = '133.58'
$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("D32").Value = '4.080'
$ws.Range("E32").Value = '  +1.69%  '

$ws.Range("D33").Value = '6.035'
$ws.Range("E33").Value = '  -0.70%  '

$ws.Range("D34").Value = '0.09201'
$ws.Range("E34").Value = '  +4.22%  '

$ws.Range("D35").Value = '12.43'
$ws.Range("E35").Value = '  -2.65%  '

$ws.Range("D36").Value = '0.02359'
$ws.Range("E36").Value = '  -2.30%  '

$ws.Range("D37").Value = '0.06347'
$ws.Range("E37").Value = '  -0.35%  '

$ws.Range("D38").Value = '0.6692'
$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("D39").Value = '5.269'
$ws.Range("E39").Value = '  -1.17%  '

$ws.Range("E40").Value = '  +0.77%  '

$ws.Range("D41").Value = '1.516'
$ws.Range("E41").Value = '  +1.15%  '

$ws.Range("D42").Value = '1.224'
$ws.Range("E42").Value = '  -1.46%  '

$ws.Range("D43").Value = '8.146'
$ws.Range("E43").Value = '  -1.78%  '

$ws.Range("D44").Value = '14.35'
$ws.Range("E44").Value = '  +1.87%  '

$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  -0.20%  '

$ws.Range("D46").Value = '0.6184'
$ws.Range("E46").Value = '  -1.38%  '

$ws.Range("D47").Value = '3.875'
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '128.75'
$ws.Range("E48").Value = '  -3.09%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '2.067'
$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").Value = '1.188'
$ws.Range("E50").Value = '  -0.45%  '

$ws.Range("D51").Value = '0.07137'
$ws.Range("E51").Value = '  -4.74%  '
